$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain numeric-looking strings (e.g. "1.01") as TEXT in the source data.
# Excel would otherwise auto-convert such literals to numbers, so we temporarily force the cell
# to Text format before assigning, then clear that temporary formatting afterwards so the cell
# keeps the workbook's default (General) style, matching the original file.

$ws.Range("D2").Value = '43.841.03'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '2.236.14'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.04'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.72'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.07'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0806'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.11'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").Value = '2.579.88'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '2.324.86'
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.824'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.40'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.75%  '
$ws.Range("D18").Value = '43.971.44'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '0.0₃0961'
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.32'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.98'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -8.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.43'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.82'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.93'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.82'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +5.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.81'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.85%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.83'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.47'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0791'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.61'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("E35").Value = '  -10.07%  '
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.108'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("E38").Value = '  -8.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.44'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.97%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0297'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.03'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.67%  '
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("D44").Value = '1.713.41'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '82.03'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.91'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '98.86'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.59'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.02'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.30'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.91%  '
